$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 398, shifting existing rows 398:425 down to 399:426
$ws.Rows.Item(398).Insert()

$ws.Cells.Item(398, 1).Value = 5
$ws.Cells.Item(398, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(398, 3).Value = "Maule"
$ws.Cells.Item(398, 4).Value = 44826
$ws.Cells.Item(398, 4).NumberFormat = $ws.Cells.Item(399, 4).NumberFormat
$ws.Cells.Item(398, 5).Value = 7
$ws.Cells.Item(398, 6).Value = 100112032
$ws.Cells.Item(398, 7).Value = "Zapallo italiano"
$ws.Cells.Item(398, 8).Value = "Sin especificar"
$ws.Cells.Item(398, 9).Value = "Primera"
$ws.Cells.Item(398, 10).Value = 300
$ws.Cells.Item(398, 11).Value = 15000
$ws.Cells.Item(398, 12).Value = 15000
$ws.Cells.Item(398, 13).Value = 15000
$ws.Cells.Item(398, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(398, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(398, 16).Value = 300
$ws.Cells.Item(398, 17).Value = 50
$ws.Cells.Item(398, 18).Value = "Hortaliza"
